$d = $word.ActiveDocument

# Section 1: "i link alla homepage" -> "il link alla homepage"
$d.Content.Find.Execute("i link alla homepage", $true, $false, $false, $false, $false, $true, 1, $false, "il link alla homepage", 2) | Out-Null

# Section 2: Reorder "Molti tasti rapidi non sono definiti per navigare efficientemente"
#         -> "Molti tasti rapidi per navigare efficientemente non sono definiti"
$d.Content.Find.Execute("Molti tasti rapidi non sono definiti per navigare efficientemente", $true, $false, $false, $false, $false, $true, 1, $false, "Molti tasti rapidi per navigare efficientemente non sono definiti", 2) | Out-Null

# Section 3: merge runs - text stays same, no find/replace needed for content but we want run merges.
